# Add two new data columns (AC, AD) with a header string in row 1 and two
# probability values in rows 2-3, extending the existing table that runs
# from column A to column AB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

$ws.Range("AC2").Value = 0.59485530546623799
$ws.Range("AD2").Value = 0.49838187702265369

$ws.Range("AC3").Value = 0.56211180124223603
$ws.Range("AD3").Value = 0.70648464163822533

# Extend the selected/active range to cover the newly added columns, same as
# the author re-selecting the whole table before saving.
$excel.Goto($ws.Range("A1:AD3"))
